# Update the "top-level-description" (column B) text for a number of
# modules on the "Specification" sheet. The previous descriptions were
# multi-line, verbose restatements of the module name; they are replaced
# with shorter, single-line, user-facing descriptions.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Specification")

$updates = @{
    21  = "Details of any changes the proposed development would make to existing access arrangements or public rights of way"
    28  = "Name and contact information if an agent is being used."
    32  = "Name and contact information if an agent is being used."
    40  = "Telephone number and email address of the applicant."
    44  = "Name and contact information for the parties making the application."
    50  = "Details of potential impacts to the biodiversity of the site, or any noteable archaeological or geological features."
    53  = "How any natural habitats on the development site will be improved by the proposed works."
    68  = "Checking whether all the requirements of the form have been met, such as proof of payment or supporting documentation."
    69  = "Details of any conflict of interest that may exist between the applicant and planning authority."
    72  = "Signed and dated verification of the application's accuracy."
    75  = "How the proposed development will impact employment, including existing and proposed employee numbers"
    82  = "How the site is currently being used."
    92  = "Results of any flood risk assessments made for the development site"
    98  = "How waste water will leave the property as part of the proposed development"
    102 = "What hazardous substances may be used as part of the development"
    108 = "Proposed operating hours if the proposed development is intended for non-residential use."
    115 = "What materials are being used for the proposed development"
    123 = "Details of changes to non-residential floorspace in the proposed development."
    134 = "Who will be affected by the proposal and whether they have been notified, such as agricultural tenants"
    149 = "Details of pre-application advice received from the planning authority"
    154 = "How waste will be managed on the site"
    164 = "What development, works or change of use is proposed"
    171 = "Details of the residential units that make up both the current and proposed development."
    187 = "How big the site is including relevant measurements"
    189 = "Where the proposed development will be built."
    198 = "Information to help the planning authority arrange a site visit"
    204 = "Details of any liquid waste produced by industial processes on the proposed site, and how it will be diposed of."
    206 = "Details of trees and/or hedges that will be affected by the proposed development"
    207 = "Details of current parking facilities at the site and any changes that would be made by the proposed development."
    212 = "Any waste storage or recycling arrangements are in place, such as waste storage areas"
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 2).Value = $updates[$row]
}
